$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2153.4443
$ws.Range("I15").Value = 2153.4443
$ws.Range("K15").Value = 6460.3329
$ws.Range("M15").Value = -6291.3329
$ws.Range("H69").Value = 9650.879999999999
$ws.Range("I69").Value = 6908
$ws.Range("K69").Value = 20724
$ws.Range("M69").Value = -19850
$ws.Range("H72").Value = 9650.879999999999
$ws.Range("I72").Value = 6908
$ws.Range("K72").Value = 62172
$ws.Range("M72").Value = -57804
$ws.Range("H103").Value = 3591.8
$ws.Range("I103").Value = 5277
$ws.Range("J103").Value = 2468.3333
$ws.Range("K103").Value = 15831
$ws.Range("L103").Value = 7404.999899999999
$ws.Range("M103").Value = -15245
$ws.Range("N103").Value = -8576.999899999999
$ws.Range("H109").Value = 99900
$ws.Range("J109").Value = 99900
$ws.Range("L109").Value = 99900
$ws.Range("N109").Value = -102674
$ws.Range("H112").Value = 4148.4375
$ws.Range("J112").Value = 3944.2307
$ws.Range("L112").Value = 11832.6921
$ws.Range("N112").Value = -14048.6921
$ws.Range("H124").Value = 102381.5
$ws.Range("J124").Value = 102381.5
$ws.Range("L124").Value = 102381.5
$ws.Range("N124").Value = -112201.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3556.6086
$ws.Range("J45").Value = 4950.3
$ws.Range("L45").Value = 4950.3
$ws.Range("N45").Value = -5704.3
$ws.Range("H53").Value = 9599.299999999999
$ws.Range("I53").Value = 5374.25
$ws.Range("J53").Value = 26499.5
$ws.Range("K53").Value = 5374.25
$ws.Range("L53").Value = 26499.5
$ws.Range("M53").Value = -4692.25
$ws.Range("N53").Value = -27863.5
$ws.Range("H97").Value = 12456.1
$ws.Range("I97").Value = 20776.4
$ws.Range("J97").Value = 4135.8
$ws.Range("K97").Value = 20776.4
$ws.Range("L97").Value = 4135.8
$ws.Range("M97").Value = -20280.4
$ws.Range("N97").Value = -5127.8
$ws.Range("H102").Value = 2122.111
$ws.Range("I102").Value = 2375
$ws.Range("J102").Value = 1237
$ws.Range("K102").Value = 2375
$ws.Range("L102").Value = 1237
$ws.Range("M102").Value = -753
$ws.Range("N102").Value = -4481

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 35665.668
$ws.Range("I26").Value = 35665.668
$ws.Range("K26").Value = 35665.668
$ws.Range("M26").Value = -35373.668
$ws.Range("H50").Value = 74999
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()
$ws.Range("H86").Value = 3197.7878
$ws.Range("I86").Value = 1696.4117
$ws.Range("K86").Value = 1696.4117
$ws.Range("M86").Value = -573.4117000000001
$ws.Range("H89").Value = 3197.7878
$ws.Range("I89").Value = 1696.4117
$ws.Range("K89").Value = 8482.058500000001
$ws.Range("M89").Value = -2866.058500000001
$ws.Range("H94").Value = 1631.4814
$ws.Range("I94").Value = 1564.4546
$ws.Range("K94").Value = 1564.4546
$ws.Range("M94").Value = -1113.4546
$ws.Range("H99").Value = 107643.3
$ws.Range("I99").Value = 149497.58
$ws.Range("K99").Value = 149497.58
$ws.Range("M99").Value = -147999.58
$ws.Range("H104").Value = 69999
$ws.Range("J104").Value = 69999
$ws.Range("L104").Value = 69999
$ws.Range("N104").Value = -76987
$ws.Range("H105").Value = 3137.8845
$ws.Range("I105").Value = 3210.8333
$ws.Range("J105").Value = 2262.5
$ws.Range("K105").Value = 3210.8333
$ws.Range("L105").Value = 2262.5
$ws.Range("M105").Value = -1463.8333
$ws.Range("N105").Value = -5756.5
$ws.Range("H115").Value = 84999
$ws.Range("J115").Value = 84999
$ws.Range("L115").Value = 84999
$ws.Range("N115").Value = -88133

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 304
$ws.Range("I7").Value = 115.8
$ws.Range("J7").Value = 438.42856
$ws.Range("K7").Value = 115.8
$ws.Range("L7").Value = 438.42856
$ws.Range("M7").Value = -2.799999999999997
$ws.Range("N7").Value = -664.4285600000001
$ws.Range("H36").Value = 12216.2
$ws.Range("I36").Value = 14024
$ws.Range("J36").Value = 11011
$ws.Range("K36").Value = 14024
$ws.Range("L36").Value = 11011
$ws.Range("M36").Value = -13636
$ws.Range("N36").Value = -11787
$ws.Range("H40").Value = 12216.2
$ws.Range("I40").Value = 14024
$ws.Range("J40").Value = 11011
$ws.Range("K40").Value = 14024
$ws.Range("L40").Value = 11011
$ws.Range("M40").Value = -13864
$ws.Range("N40").Value = -11331
$ws.Range("H53").Value = 43144.5
$ws.Range("I53").Value = 11289
$ws.Range("J53").Value = 75000
$ws.Range("K53").Value = 11289
$ws.Range("L53").Value = 75000
$ws.Range("M53").Value = -10682
$ws.Range("N53").Value = -76214
$ws.Range("H62").Value = 4358.4287
$ws.Range("J62").Value = 3402.75
$ws.Range("L62").Value = 3402.75
$ws.Range("N62").Value = -4650.75
$ws.Range("H65").Value = 4358.4287
$ws.Range("J65").Value = 3402.75
$ws.Range("L65").Value = 17013.75
$ws.Range("N65").Value = -23253.75
$ws.Range("H105").Value = 2018.6666
$ws.Range("I105").Value = 2018.6666
$ws.Range("K105").Value = 2018.6666
$ws.Range("M105").Value = -271.6666
$ws.Range("H137").Value = 99999
$ws.Range("J137").Value = 99999
$ws.Range("L137").Value = 99999
$ws.Range("N137").Value = -110199
$ws.Range("H139").Value = 79999
$ws.Range("J139").Value = 79999
$ws.Range("L139").Value = 79999
$ws.Range("N139").Value = -90279

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 448
$ws.Range("I86").Value = 395.66666
$ws.Range("K86").Value = 1186.99998
$ws.Range("M86").Value = -0.9999800000000505
$ws.Range("H89").Value = 448
$ws.Range("I89").Value = 395.66666
$ws.Range("K89").Value = 3560.99994
$ws.Range("M89").Value = 2367.00006
$ws.Range("H98").Value = 1979.1666
$ws.Range("I98").Value = 1900
$ws.Range("J98").Value = 1995
$ws.Range("K98").Value = 5700
$ws.Range("L98").Value = 5985
$ws.Range("M98").Value = -4202
$ws.Range("N98").Value = -8981
$ws.Range("H134").Value = 798.8
$ws.Range("I134").Value = 798.8
$ws.Range("K134").Value = 2396.4
$ws.Range("M134").Value = 2673.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 10499.5
$ws.Range("I80").Value = 2997
$ws.Range("J80").Value = 13000.333
$ws.Range("K80").Value = 2997
$ws.Range("L80").Value = 13000.333
$ws.Range("M80").Value = -1999
$ws.Range("N80").Value = -14996.333
$ws.Range("H83").Value = 10499.5
$ws.Range("I83").Value = 2997
$ws.Range("J83").Value = 13000.333
$ws.Range("K83").Value = 14985
$ws.Range("L83").Value = 65001.665
$ws.Range("M83").Value = -9993
$ws.Range("N83").Value = -74985.66500000001
$ws.Range("H132").Value = 54969.633
$ws.Range("I132").Value = 54969.633
$ws.Range("K132").Value = 164908.899
$ws.Range("M132").Value = -162378.899

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3938.8
$ws.Range("J7").Value = 5895.5
$ws.Range("L7").Value = 5895.5
$ws.Range("N7").Value = -6119.5
$ws.Range("H82").Value = 3004.923
$ws.Range("I82").Value = 2698
$ws.Range("K82").Value = 2698
$ws.Range("M82").Value = -2337
$ws.Range("H85").Value = 3004.923
$ws.Range("I85").Value = 2698
$ws.Range("K85").Value = 2698
$ws.Range("M85").Value = -1450
$ws.Range("H93").Value = 3759.8
$ws.Range("J93").Value = 2999
$ws.Range("L93").Value = 2999
$ws.Range("N93").Value = -5495
$ws.Range("H126").Value = 3938.8
$ws.Range("J126").Value = 5895.5
$ws.Range("L126").Value = 17686.5
$ws.Range("N126").Value = -22626.5
$ws.Range("H132").Value = 89563.5
$ws.Range("I132").Value = 152698.75
$ws.Range("J132").Value = 5383.1665
$ws.Range("K132").Value = 458096.25
$ws.Range("L132").Value = 16149.4995
$ws.Range("M132").Value = -455566.25
$ws.Range("N132").Value = -21209.4995
$ws.Range("H133").Value = 93318.336
$ws.Range("J133").Value = 93318.336
$ws.Range("L133").Value = 93318.336
$ws.Range("N133").Value = -98378.336
$ws.Range("H136").Value = 3045.8235
$ws.Range("I136").Value = 2052
$ws.Range("K136").Value = 6156
$ws.Range("M136").Value = -3606

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1555.6923
$ws.Range("I81").Value = 1063
$ws.Range("K81").Value = 2126
$ws.Range("M81").Value = -1065
$ws.Range("H84").Value = 1555.6923
$ws.Range("I84").Value = 1063
$ws.Range("K84").Value = 10630
$ws.Range("M84").Value = -5326
$ws.Range("H124").Value = 67210.39999999999
$ws.Range("J124").Value = 67210.39999999999
$ws.Range("L124").Value = 67210.39999999999
$ws.Range("N124").Value = -77030.39999999999
